# Updated cryptos list (prices + 1h volume change) per the authored diff.
# Column D (Price) holds text-formatted numbers (e.g. "56.891.43", trailing
# zeros like "4.20") in the source data, so every D-column write is guarded
# with a Text number-format (then ClearFormats to drop the now-unneeded
# style) to stop Excel from silently re-typing the cell as a float and
# mangling the literal text (dropping trailing zeros / FP rounding noise).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.891.43"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.347.60"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "515.54"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.17"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.535"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.39"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.74%  "
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.342"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.03"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.765.57"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.837.90"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.356.25"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.46"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "327.49"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.20"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.71"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.13"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.166"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.69"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +12.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.996"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.32"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +8.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.46"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0731"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("E30").Value = "  +0.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.20"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.51"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.37%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("E34").Value = "  +2.98%  "
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("E36").Value = "  +0.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.895"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.57"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.57%  "
$ws.Range("E39").Value = "  +3.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "150.52"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +8.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.379"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("E42").Value = "  +1.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "283.49"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.58%  "
$ws.Range("E44").Value = "  +4.28%  "
$ws.Range("E45").Value = "  +0.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0502"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.559"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.43"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +8.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0217"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.31"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.39%  "
$ws.Range("E51").Value = "  +1.07%  "
